$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14: Two Sum
$ws.Range("A14").Value = "Two Sum"
$ws.Range("B14").Value = "Return idexes of the target sum in an array"
$ws.Range("C14").Value = "Create a complement hashmap to keep track of complements(target - currentNum). If complement exist for currentNum in map, we found our target. "
$ws.Range("D14").Value = "https://leetcode.com/problems/two-sum/"
$ws.Hyperlinks.Add($ws.Range("D14"), "https://leetcode.com/problems/two-sum/")

# Row 15: Two Sum II - Input Array Is Sorted
$ws.Range("A15").Value = "Two Sum II - Input Array Is Sorted"
$ws.Range("B15").Value = "Return idexes of the target sum in an array"
$ws.Range("C15").Value = "Use two pointers left and right to calculate current sum. If current sum > target, decrease right pointer else increase left pointer"
$ws.Range("D15").Value = "https://leetcode.com/problems/two-sum-ii-input-array-is-sorted/"
$ws.Hyperlinks.Add($ws.Range("D15"), "https://leetcode.com/problems/two-sum-ii-input-array-is-sorted/")

# Match styles used by other rows: column A uses the "Good" look (same as
# rows 2-5 & 12), columns B-D reuse row 13's formatting.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null

$ws.Range("B13:D13").Copy() | Out-Null
$ws.Range("B14:D14").PasteSpecial(-4122) | Out-Null
$ws.Range("B15:D15").PasteSpecial(-4122) | Out-Null

$ws.Range("D16").Select()
